# Update the yearly vehicle-uptake projections on sheet "a" (rows 1-7).
# Columns I, J, K get new figures; columns L through Q are reset to the
# #N/A error value (matching the pattern already used from column R on).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("a")

$rows = @(1, 2, 3, 4, 5, 6, 7)
$newIJK = @{
    1 = @(109, 216, 312)
    2 = @(109, 216, 312)
    3 = @(159, 304, 439.5)
    4 = @(159, 304, 439.5)
    5 = @(99.5, 207, 311.5)
    6 = @(99.5, 207, 311.5)
    7 = @(503, 975, 1447)
}

foreach ($r in $rows) {
    $vals = $newIJK[$r]
    $ws.Range("I$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]

    foreach ($col in @("L", "M", "N", "O", "P", "Q")) {
        $ws.Range("$col$r").Value = "#N/A"
    }
}

# The little summary table (rows 10, 11, 14, 15) stores its "Emisiones CO2"
# (col C) and "Numero de Spark requeridos" (col D) as plain numbers rather
# than live formulas, so they need to be refreshed by hand to stay in sync
# with the updated vehicle counts above (mirrors Hoja2!E14:H17, which is
# `=+a!G*a!E*a!I` / `=+(a!G+a!H)*a!$E$7` / `=E/G`, floored for column D).
$summary = @(
    @{ Row = 10; G = "G1"; E = "E1"; H = "H1"; I = "I1" }
    @{ Row = 11; G = "G2"; E = "E2"; H = "H2"; I = "I2" }
    @{ Row = 14; G = "G5"; E = "E5"; H = "H5"; I = "I5" }
    @{ Row = 15; G = "G6"; E = "E6"; H = "H6"; I = "I6" }
)

$E7 = $ws.Range("E7").Value2

foreach ($item in $summary) {
    $g = $ws.Range($item.G).Value2
    $e = $ws.Range($item.E).Value2
    $h = $ws.Range($item.H).Value2
    $i = $ws.Range($item.I).Value2

    $emissions = $g * $e * $i
    $vehicles = ($g + $h) * $E7
    $ratio = $emissions / $vehicles

    $ws.Range("C$($item.Row)").Value = $emissions / 1000
    $ws.Range("D$($item.Row)").Value = [Math]::Floor($ratio)
}
